$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row stays pattern/responses/context_set (C1 string index changes but text is the same)
$ws.Range("A1").Value = "pattern"
$ws.Range("B1").Value = "responses"
$ws.Range("C1").Value = "context_set"

$website1 = "You can find the current course website at http://www.sba.oakland.edu/faculty/isken/courses/mis6900_s24/, You can find the permalink at http://www.sba.oakland.edu/faculty/isken/courses/aap"
$website2 = "You can find the current course website at http://www.sba.oakland.edu/faculty/isken/courses/mis5460_f23/index.html, You can find the permalink at http://www.sba.oakland.edu/faculty/isken/courses/ba"
$website3 = "You can find the current course website at http://www.sba.oakland.edu/faculty/isken/courses/mis5470_f23/index.html, You can find the permalink at http://www.sba.oakland.edu/faculty/isken/courses/pcda"

$ws.Range("A2").Value = "mis4900"
$ws.Range("B2").Value = $website1

$ws.Range("A3").Value = "mis4460"
$ws.Range("B3").Value = $website2

$ws.Range("A4").Value = "mis4470"
$ws.Range("B4").Value = $website3

$ws.Range("A5").Value = "mis6900"
$ws.Range("B5").Value = $website1

$ws.Range("A6").Value = "mis5460"
$ws.Range("B6").Value = $website2

$ws.Range("A7").Value = "mis5470"
$ws.Range("B7").Value = $website3

$ws.Range("A8").Value = "quit"
$ws.Range("B8").Value = "Good bye!"

$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null

$ws.Range("A1").Select() | Out-Null
